# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Source values come from the coinranking.com feed; D = Price, E = Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.472.78'
$ws.Range('E2').Value = '  -0.99%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.104.42'
$ws.Range('E3').Value = '  -0.34%  '

# Row 4: TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.39'
$ws.Range('E5').Value = '  -0.32%  '

# Row 6: USDC
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.12%  '

# Row 7: XRP
$ws.Range('E7').Value = '  -1.56%  '

# Row 8: Cardano
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4505'
$ws.Range('E8').Value = '  +2.46%  '

# Row 9: OKB
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.79'
$ws.Range('E9').Value = '  +16.80%  '

# Row 10: Dogecoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08931'
$ws.Range('E10').Value = '  -1.00%  '

# Row 11: Polygon
$ws.Range('E11').Value = '  -1.81%  '

# Row 12: Solana
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.32'
$ws.Range('E12').Value = '  -2.75%  '

# Row 13: WrappedEther
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.089.69'
$ws.Range('E13').Value = '  -0.74%  '

# Row 14: Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.742'
$ws.Range('E14').Value = '  -0.19%  '

# Row 15: Chainlink
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.750'
$ws.Range('E15').Value = '  -0.69%  '

# Row 16: Litecoin
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.31'
$ws.Range('E16').Value = '  -1.17%  '

# Row 17: BinanceUSD
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  +0.10%  '

# Row 18: ShibaInu
$ws.Range('E18').Value = '  -0.28%  '

# Row 19: TRON
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06610'
$ws.Range('E19').Value = '  -0.78%  '

# Row 20: Avalanche
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.27'
$ws.Range('E20').Value = '  +0.74%  '

# Row 21: Dai
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.07%  '

# Row 22: Uniswap
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.288'
$ws.Range('E22').Value = '  -1.07%  '

# Row 23: WrappedBTC
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.519.96'
$ws.Range('E23').Value = '  -1.02%  '

# Row 24: Cosmos
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.31'
$ws.Range('E24').Value = '  -0.49%  '

# Row 25: Toncoin
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.323'
$ws.Range('E25').Value = '  +2.89%  '

# Row 26: WrappedliquidstakedEther2.0
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.340.19'
$ws.Range('E26').Value = '  -0.60%  '

# Row 27: EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.29'
$ws.Range('E27').Value = '  -2.11%  '

# Row 28: LidoDAOToken
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.585'
$ws.Range('E28').Value = '  +0.44%  '

# Row 29: Monero
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.83'
$ws.Range('E29').Value = '  +0.86%  '

# Row 30: BitcoinCash
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.78'
$ws.Range('E30').Value = '  -0.08%  '

# Row 31: ImmutableX
$ws.Range('E31').Value = '  +2.31%  '

# Row 32: Stellar
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1074'
$ws.Range('E32').Value = '  -0.41%  '

# Row 33: ARBITRUM
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.672'
$ws.Range('E33').Value = '  +8.40%  '

# Row 34: Filecoin
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.146'
$ws.Range('E34').Value = '  -1.29%  '

# Row 35: HuobiToken
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.944'
$ws.Range('E35').Value = '  -1.84%  '

# Row 36: FraxShare
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.43'
$ws.Range('E36').Value = '  +9.52%  '

# Row 37: VeChain
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02568'
$ws.Range('E37').Value = '  -1.43%  '

# Row 39: Aptos
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.77'
$ws.Range('E39').Value = '  -0.12%  '

# Row 40: InternetComputer(DFINITY)
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.483'
$ws.Range('E40').Value = '  -0.93%  '

# Row 41: Algorand
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2266'
$ws.Range('E41').Value = '  -0.38%  '

# Row 42: TheSandbox
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6915'
$ws.Range('E42').Value = '  +0.74%  '

# Row 43: TrustWalletToken
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.257'
$ws.Range('E43').Value = '  +0.35%  '

# Row 44: Frax
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.11%  '

# Row 45: EnergySwap
$ws.Range('E45').Value = '  -1.25%  '

# Row 46: ranking reshuffled, Decentraland -> NEARProtocol
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.298'
$ws.Range('E46').Value = '  +2.96%  '

# Row 47: ranking reshuffled, NEARProtocol -> Decentraland
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6361'
$ws.Range('E47').Value = '  -1.36%  '

# Row 48: PancakeSwap
$ws.Range('E48').Value = '  -0.82%  '

# Row 49: EOS
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.244'
$ws.Range('E49').Value = '  -2.47%  '

# Row 50: WEMIXTOKEN
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.225'
$ws.Range('E50').Value = '  +5.80%  '

# Row 51: Aave
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '82.71'
$ws.Range('E51').Value = '  +0.39%  '
